# Journal update: add "sprint review" / "sprint 2" week entries (week of 2025-12-18)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal_NOM_PRENOM")

# --- Week 4 block (rows 20-24), first daily entry gets the date ---
$ws.Range("A20").Value = 46009
$ws.Range("B20").Value = "Preparation sprint review et correction de problème"
$ws.Range("D20").Value = 2

# --- Week 4 reflection box (merged B26:D26) ---
$ws.Range("B26").Value = "Petit problème avec les jdk qui se supprime mais rien d'alarmantl"

# --- Remaining week 4 daily entries ---
$ws.Range("B21").Value = "sprint review"
$ws.Range("D21").Value = 1

$ws.Range("B22").Value = "commencement sprint 2"
$ws.Range("D22").Value = 0.5

# --- Selection / view bookkeeping to match the saved workbook state ---
$ws.Range("B22:C22").Select()
